# Xong sign in và xử lý ngày
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-wide font change: Calibri -> Arial (update the base/Normal style
# so existing cells keep using the default style index, matching the saved
# file which carries no explicit per-cell "s" attribute)
$wb.Styles.Item("Normal").Font.Name = "Arial"

# New columns: Ngày / Tổng thời gian fix / Kinh nghiệm rút ra
$ws.Range("I1").Value = "Ngày"
$ws.Range("I2").Value = "28/10/2025"

$ws.Range("J1").Value = "Tổng thời gian fix"
$ws.Range("J2").Value = "4 giờ"

$ws.Range("K1").Value = "Kinh nghiệm rút ra"
$ws.Range("K2").Value = "Kiểm tra kĩ CSS selector"

# Match the saved view state from the diff (active cell selection)
$ws.Range("I10").Select()
